$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update metadata for the "municipio-nombre" column (F) to reflect the
# newly curated dimensions.
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("F3").Value = "dim"
$ws.Range("F4").Value = "URI-Municipio"
